# Apply the data-repair edit described by the commit diff for LOM3115.xlsx
# The edit adds a new shared string "1033242 - Fábio Herbst Florenzano" and
# shuffles the shared-string table order, which (together with the
# unchanged cell->sharedstring references for rows 1-12) results in these
# effective textual changes on the single worksheet:
#
#   B10/C10 : "Projetos pré-estruturados ..." -> "1033242 - Fábio Herbst Florenzano"
#   B13/C13 : "30 h"        -> "Semestral"
#   B15/C15 : "Semestral"   -> "01/01/2023"
#   B18/C18 : "01/01/2023"  -> "1033242 - Fábio Herbst Florenzano"
#
# (B19/C19, B20/C20, B21/C21, A13, A14, A15, A18, etc. keep the same text,
#  only the underlying shared-string index shifts, so no visible change.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "1033242 - Fábio Herbst Florenzano"
$ws.Range("C10").Value = "1033242 - Fábio Herbst Florenzano"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B15").NumberFormat = "@"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("B15").Value = "01/01/2023"
$ws.Range("C15").Value = "01/01/2023"

$ws.Range("B18").Value = "1033242 - Fábio Herbst Florenzano"
$ws.Range("C18").Value = "1033242 - Fábio Herbst Florenzano"
